$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the existing header row (bold, border, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every player row (same record for the whole team)
$ws.Range("AD2:AD46").Value = 83
$ws.Range("AE2:AE46").Value = 79
$ws.Range("AF2:AF46").Value = 0
